$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-converted to a number by Excel (losing the original text formatting).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values.
$ws.Range("D2").Value = '29.220.50'
$ws.Range("D3").Value = '1.862.26'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '0.7041'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").Value = '237.87'
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '0.08194'
$ws.Range("E8").Value = '  +9.85%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").Value = '23.29'
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("D11").Value = '0.08161'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '1.839.41'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").Value = '5.167'
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("D14").Value = '0.7087'
$ws.Range("E14").Value = '  -2.25%  '
$ws.Range("D15").Value = '89.19'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '29.224.30'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '0.000007883'
$ws.Range("E17").Value = '  +3.44%  '
$ws.Range("D18").Value = '5.783'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("E19").Value = '  +2.25%  '
$ws.Range("D20").Value = '236.07'
$ws.Range("E20").Value = '  -0.83%  '
$ws.Range("D21").Value = '0.9990'
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").Value = '2.107.41'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").Value = '0.9999'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '7.403'
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("D25").Value = '162.04'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("D26").Value = '8.955'
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").Value = '0.1441'
$ws.Range("E27").Value = '  -0.65%  '
$ws.Range("D28").Value = '18.07'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '1.963'
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("D30").Value = '1.427'
$ws.Range("E30").Value = '  +2.12%  '
$ws.Range("D31").Value = '1.485'
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("E32").Value = '  -3.47%  '
$ws.Range("D33").Value = '4.059'
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").Value = '1.169'
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("D36").Value = '0.7068'
$ws.Range("E36").Value = '  +1.03%  '
$ws.Range("D37").Value = '0.9983'
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("D38").Value = '2.669'
$ws.Range("E38").Value = '  +0.62%  '
$ws.Range("D39").Value = '0.01849'
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("D40").Value = '2.730'
$ws.Range("E40").Value = '  +1.99%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.9246'
$ws.Range("E41").Value = '  -1.38%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.142.04'
$ws.Range("E42").Value = '  +6.00%  '
$ws.Range("D43").Value = '0.4276'
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").Value = '5.871'
$ws.Range("E44").Value = '  -2.41%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '0.9989'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").Value = '102.91'
$ws.Range("D48").Value = '1.775'
$ws.Range("E48").Value = '  +1.89%  '
$ws.Range("D49").Value = '2.000.51'
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("D51").Value = '6.956'
$ws.Range("E51").Value = '  -1.08%  '
